$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("E2").Value = 0.000260678
$ws.Range("F2").Value = 0.017472907
$ws.Range("G2").Value = 0.0004503267894

$ws.Range("E3").Value = 0.002693689
$ws.Range("F3").Value = 0.007200156
$ws.Range("G3").Value = 0.003091181772136223

$ws.Range("E4").Value = 0.005545441
$ws.Range("F4").Value = 0.01091392
$ws.Range("G4").Value = 0.006741560248313091
